# Refresh the crypto price/ranking table (GitHub Actions "Updated symbol list" commit).
# Most edits are text-valued "Price" (column D) updates; a handful of rows also
# had their Coin/Link/Volume(1h) (columns B/C/E) replaced because the underlying
# coin ranking shifted by one position for two short runs of rows (9-17 and 42-43).
#
# All of these cells are stored as *text* (not numbers) in the source workbook,
# so assigning a numeric-looking value needs to be forced to text (leading
# apostrophe) to avoid Excel's automatic number coercion - otherwise values like
# "0.1300" or "0.04110" would lose their trailing zeros, and "247.71" would be
# stored as a float instead of a string. ClearFormats() afterwards drops the
# "quote prefix" text style Excel applies automatically, so the cell keeps the
# same (default) style it had before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.ClearFormats()
}

# --- Rows 2-8: only the Price (D) column changed ---
Set-TextCell 2 4 '247.71'
Set-TextCell 3 4 '21.66'
Set-TextCell 4 4 '5.495'
Set-TextCell 5 4 '0.05684'
Set-TextCell 6 4 '3.387'
Set-TextCell 7 4 '0.8081'
Set-TextCell 8 4 '1.035'

# --- Rows 9-17: Coin/Link/Price/Volume(1h) all shifted down one rank ---
Set-TextCell 9 2 'One'
Set-TextCell 9 3 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell 9 4 '0.01160'
Set-TextCell 9 5 '8OneONEBestin24h'

Set-TextCell 10 2 'WazirX'
Set-TextCell 10 3 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 10 4 '0.1519'
Set-TextCell 10 5 '9WazirXWRX'

Set-TextCell 11 2 'MandalaExchangeToken'
Set-TextCell 11 3 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 11 4 '0.07908'
Set-TextCell 11 5 '10MandalaExchangeTokenMDX'

Set-TextCell 12 2 'LiechtensteinCryptoassetsExchange'
Set-TextCell 12 3 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 12 4 '0.03153'
Set-TextCell 12 5 '11LiechtensteinCryptoassetsExchangeLCX'

Set-TextCell 13 2 'BitrueCoin'
Set-TextCell 13 3 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 13 4 '0.03024'
Set-TextCell 13 5 '12BitrueCoinBTR'

Set-TextCell 14 2 'BitMartToken'
Set-TextCell 14 3 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 14 4 '0.09294'
Set-TextCell 14 5 '13BitMartTokenBMX'

Set-TextCell 15 2 'MCDex'
Set-TextCell 15 3 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextCell 15 4 '3.466'
Set-TextCell 15 5 '14MCDexMCB'

Set-TextCell 16 2 'BitForexToken'
Set-TextCell 16 3 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 16 4 '0.001657'
Set-TextCell 16 5 '15BitForexTokenBF'

Set-TextCell 17 2 'CoinExToken'
Set-TextCell 17 3 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextCell 17 4 '0.04722'
Set-TextCell 17 5 '16CoinExTokenCET'

# --- Rows 18-27: only Price (D) changed (row 19 also drops "Bestin24h" from E) ---
Set-TextCell 18 4 '0.006360'

Set-TextCell 19 4 '0.005034'
Set-TextCell 19 5 '18HotbitTokenHTB'

Set-TextCell 20 4 '0.001044'
Set-TextCell 21 4 '0.0001502'
Set-TextCell 22 4 '0.0003204'
Set-TextCell 23 4 '3.766'
Set-TextCell 24 4 '6.425'
Set-TextCell 25 4 '2.152'
Set-TextCell 26 4 '0.3314'
Set-TextCell 27 4 '0.1300'

# --- Rows 40-41: only Price (D) changed ---
Set-TextCell 40 4 '0.04110'
Set-TextCell 41 4 '0.006964'

# --- Rows 42-43: Coin/Link/Price/Volume(1h) swapped ranking ---
Set-TextCell 42 2 'CEJI'
Set-TextCell 42 3 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell 42 4 '0.003505'
Set-TextCell 42 5 '41CEJICEJI'

Set-TextCell 43 2 'BKEXToken'
Set-TextCell 43 3 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell 43 4 '0.1044'
Set-TextCell 43 5 '42BKEXTokenBKK'

# --- Rows 44-50: only Price (D) changed ---
Set-TextCell 44 4 '0.007766'
Set-TextCell 45 4 '0.00005892'
Set-TextCell 48 4 '0.6834'
Set-TextCell 49 4 '0.008497'
Set-TextCell 50 4 '0.00002103'
